$d = $word.ActiveDocument

# Near the end of the document, three paragraphs are being removed (the
# "scraped site" boilerplate that Jekyll appended on build):
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#     pages. Original theme under Creative Commons Attribution"
#   - the blank paragraph right after that one
# Everything from the start of the "Ver no Jupiter..." paragraph through
# the end of the blank paragraph that follows the copyright notice gets
# deleted, leaving the previously-existing blank paragraph (right after
# "LOB1012: Estatística (Requisito fraco)") directly adjacent to the
# paragraph carrying the trailing page break.

$startPar = $null
$endPar = $null

foreach ($p in $d.Paragraphs) {
    if ($startPar -eq $null -and $p.Range.Text -like "Ver no Jupiter*") {
        $startPar = $p
    }
    if ($p.Range.Text -like "*Powered by Jekyll*") {
        $endPar = $p
        break
    }
}

if ($startPar -ne $null -and $endPar -ne $null) {
    # Also sweep up the blank paragraph immediately following the
    # copyright-notice paragraph.
    $afterEnd = $endPar.Next()
    if ($afterEnd -ne $null) {
        $endPar = $afterEnd
    }

    $delStart = $startPar.Range.Start
    $delEnd = $endPar.Range.End
    $d.Range($delStart, $delEnd).Delete()
}
